# Applies the "#5: insurance, claim, debt, investment done" edit:
#  - Sheet "債務" (debt, sheet5): adds proper column headers and appends
#    the standard property_category/category/date/legislator_name/
#    legislator_id/source_file/index columns (H:N) to every data row.
#  - Sheet "事業投資" (investment, sheet6): same header/column completion,
#    plus fixes the "total" value on row 2 which was stored as text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 5: 債務 (debt)
# ---------------------------------------------------------------------
$debt = $wb.Worksheets.Item("債務")

# Header row (row 1) - B:G get their real column names, H:N are the
# shared trailer columns used across every sheet in the workbook.
$debt.Range("B1").Value = "species"
$debt.Range("C1").Value = "debtor"
$debt.Range("D1").Value = "owner"
$debt.Range("E1").Value = "total"
$debt.Range("F1").Value = "register_date"
$debt.Range("G1").Value = "register_reason"
$debt.Range("H1").Value = "property_category"
$debt.Range("I1").Value = "category"
$debt.Range("J1").Value = "date"
$debt.Range("K1").Value = "legislator_name"
$debt.Range("L1").Value = "legislator_id"
$debt.Range("M1").Value = "source_file"
$debt.Range("N1").Value = "index"

# Rows 2-4: fill in the new trailer columns H:N for each existing record.
for ($r = 2; $r -le 4; $r++) {
    $idx = $debt.Range("A" + $r).Value2
    $debt.Range("H" + $r).Value = "debt"
    $debt.Range("I" + $r).Value = "normal"
    # Quote-prefix so Excel stores the literal text "2011-11-17" instead
    # of auto-converting it to a date serial, then drop the resulting
    # quote-prefix number format so the cell is plain/unformatted again.
    $debt.Range("J" + $r).Value = "'2011-11-17"
    $debt.Range("J" + $r).ClearFormats()
    $debt.Range("K" + $r).Value = "呂學樟"
    $debt.Range("L" + $r).Value = 892
    $debt.Range("M" + $r).Value = "tmpf9381"
    $debt.Range("N" + $r).Value = $idx
}

# ---------------------------------------------------------------------
# Sheet 6: 事業投資 (business investment)
# ---------------------------------------------------------------------
$invest = $wb.Worksheets.Item("事業投資")

$invest.Range("B1").Value = "owner"
$invest.Range("C1").Value = "company"
$invest.Range("D1").Value = "address"
$invest.Range("E1").Value = "total"
$invest.Range("F1").Value = "register_date"
$invest.Range("G1").Value = "register_reason"
$invest.Range("H1").Value = "property_category"
$invest.Range("I1").Value = "category"
$invest.Range("J1").Value = "date"
$invest.Range("K1").Value = "legislator_name"
$invest.Range("L1").Value = "legislator_id"
$invest.Range("M1").Value = "source_file"
$invest.Range("N1").Value = "index"

# Row 2's acquire value was stored as text "13800000" - make it numeric.
$invest.Range("E2").Value = 13800000

for ($r = 2; $r -le 3; $r++) {
    $idx = $invest.Range("A" + $r).Value2
    $invest.Range("H" + $r).Value = "investment"
    $invest.Range("I" + $r).Value = "normal"
    $invest.Range("J" + $r).Value = "'2011-11-17"
    $invest.Range("J" + $r).ClearFormats()
    $invest.Range("K" + $r).Value = "呂學樟"
    $invest.Range("L" + $r).Value = 892
    $invest.Range("M" + $r).Value = "tmpf9381"
    $invest.Range("N" + $r).Value = $idx
}

Write-Output "done"
